# Added support for both Internal and External Resources
#
# The PowerPoint position/size properties (Left/Top/Width/Height) are
# stored internally as points (single-precision) while the OOXML stores
# EMU (1 pt = 12700 EMU). A plain "emu / 12700" division can land a hair
# below the intended boundary and round down once re-quantized, so we
# nudge by half an EMU (in points) to make the EMU round-trip exact.
function EMUToPt([double]$emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- Shape 1: "Picture 1" -- reposition/resize the diagram image ---
$pic = $s.Shapes.Item(1)
$pic.Left   = EMUToPt 180468
$pic.Top    = EMUToPt 136147
$pic.Width  = EMUToPt 12220344
$pic.Height = EMUToPt 6721853

# --- Shape 3: "TextBox 6" ("Entry Point. i.e Main() in Main.py ...") ---
$tb6 = $s.Shapes.Item(3)
$tb6.Left = EMUToPt 6460376
$tb6.Top  = EMUToPt 1391330
# Force a text rewrite so the two runs collapse into a single run.
$tb6.TextFrame.TextRange.Text = "x"
$tb6.TextFrame.TextRange.Text = "Entry Point. i.e Main() in Main.py (Package DemoStandAloneApp)"

# --- Shape 4: "TextBox 7" ("Copy and rename a Single 'pptx' File") ---
$tb7 = $s.Shapes.Item(4)
$tb7.Left = EMUToPt 5598147
$tb7.Top  = EMUToPt 3599294

# --- Shape 7: "TextBox 10" ("For Debugging Purpose. / NO, Will retain ... / after build ...") ---
$tb10 = $s.Shapes.Item(7)
$tb10.Left = EMUToPt 6038830
$tb10.Top  = EMUToPt 5437935
# Merge the two runs of the 2nd paragraph only; leave the other paragraphs untouched.
$para2 = $tb10.TextFrame.TextRange.Paragraphs(2, 1)
$para2.Text = "x"
$para2 = $tb10.TextFrame.TextRange.Paragraphs(2, 1)
$para2.Text = "NO, Will retain the temporary OBJ folder"

# --- Shape 8: "TextBox 11" ("Copy 'RunApp.sh' to BLD_DIR") ---
$tb11 = $s.Shapes.Item(8)
$tb11.Left = EMUToPt 5598147
$tb11.Top  = EMUToPt 2109044

# --- Shape 9: "TextBox 12" ("These are internal resources / ... / PYZ Archive file") ---
$tb12 = $s.Shapes.Item(9)
$tb12.Left = EMUToPt 9890346
$tb12.Top  = EMUToPt 3667407
